# Swap the presentation's applied colour theme from the custom "Integral"
# palette to the stock Office default palette ("Office Theme"), the same
# effect as picking the first/default theme in the PowerPoint Design tab.
#
# PowerPoint exposes the active design's 12-slot colour scheme through
# SlideMaster.Theme.ThemeColorScheme (Item 1..12 == dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink). Re-pointing every slot to the Office Theme
# RGB values rewrites the <a:clrScheme> that the slide master's theme part
# carries.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeRGB = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeRGB[$i - 1]
}
